$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows (rows 11-22) for the matrix-multiplication measurements
# Columns: A=Matrix Size M, B=Matrix Size P, C=Matrix Size N, D=Instructions,
#          E=Cycles, F=Cache reference, G=Cache miss, H=Cache hit ratio (formula), I=Time spent
$data = @(
    @(98,   461, 9,   23837684,    25093083,   28747,   14260,  13678482),
    @(128,  439, 612, 1132344495, 610098078,  301008,  71408,  399326069),
    @(460,  689, 809, 8093239977, 3391057734, 3005137, 621630, 1826409408),
    @(873,  987, 611, 16537254902,6606571601, 6412068, 1179174,3153030696),
    @(263,  542, 682, 3101765457, 1329933418, 893578,  44167,  850191702),
    @(46,   836, 710, 935857312,  474943232,  344619,  122864, 303112921),
    @(992,  403, 15,  242708708,  160322043,  31848,   16287,  101737111),
    @(733,  339, 385, 3051855365, 1353288115, 588402,  66978,  711956200),
    @(973,  97,  942, 2851377745, 1188223003, 542545,  100067, 773122899),
    @(654,  959, 638, 12589054691,5471941717, 4843469, 786366, 2593891840),
    @(47,   182, 700, 221359403,  140658675,  64140,   23956,  97192970),
    @(791,  197, 434, 2166540466, 1072472631, 338962,  63584,  689245593)
)

$startRow = 11
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]

    $ws.Cells.Item($row, 9).Value = $vals[7]
    $ws.Cells.Item($row, 9).Style = "Input"
}

# Extend the shared "cache hit ratio" formula down through row 22
$ws.Range("H4:H22").FormulaR1C1 = "=ROUND(((RC6 - RC7) / RC6) * 100, 0)"

# Replicate the style quirk on A11 from the original edit (copy-paste artifact)
$ws.Cells.Item(11, 1).Style = "Input"

$ws.Range("F27").Select()
